$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, D(open), E(close), F(high), G(low), H(shares_outstanding)
# I (fixed_ticker) is set to "STX" for every data row.
$rows = @(
    @(2, 40.14482874024634, 33.98491287231445, 40.46998615540762, 33.90663602187954, 212968277),
    @(3, 31.62416567886437, 35.67667770385742, 36.21741715755352, 31.01659284231479, 212968277),
    @(4, 29.35318208390361, 31.02717781066895, 32.4252420009045, 27.26222013152092, 212968277),
    @(5, 27.76902975202814, 23.58601379394531, 30.67545035927411, 22.80518282460153, 212968277),
    @(6, 22.68197645351937, 18.28785514831543, 23.21707866150206, 16.52517074874302, 212968277),
    @(7, 21.76709185261312, 13.97845363616943, 22.98065493415028, 13.87571820200242, 212968277),
    @(8, 16.05050441407019, 21.23493003845215, 21.89790009301792, 14.64500860651917, 212968277),
    @(9, 26.01437322869416, 23.14712524414062, 26.53385140766549, 22.62764706516929, 212968277),
    @(10, 26.37182825209555, 30.95912933349609, 32.00824181383994, 24.36960016990746, 212968277),
    @(11, 31.714543606655, 29.28191375732422, 35.41908996993069, 28.63552904239356, 212968277),
    @(12, 27.46248814059954, 23.25703048706055, 28.34450540141072, 22.49496762755375, 212968277),
    @(13, 24.19152753249695, 26.58623123168945, 29.4124117403204, 23.83915544488351, 212968277),
    @(14, 30.90808921992572, 40.28634262084961, 41.19862391758959, 30.58696720573542, 212968277),
    @(15, 43.51355341063856, 42.69490814208984, 46.24237003466516, 41.3526258721224, 212968277),
    @(16, 41.31452918978942, 39.2271842956543, 44.47536785628962, 38.7128037589177, 212968277),
    @(17, 36.0858345962994, 30.39005851745605, 36.49375554109401, 27.71591394378342, 212968277),
    @(18, 29.1354452573179, 34.02208709716797, 35.58181896176406, 28.02135483672853, 212968277),
    @(19, 37.6572079271242, 37.61828231811523, 39.92271282151572, 33.81130014226662, 212968277),
    @(20, 38.51562040713414, 36.57286071777344, 39.16320697025437, 36.03583727333923, 212968277),
    @(21, 43.16144301395668, 46.38256454467773, 46.47048649051457, 40.85949853281612, 212968277),
    @(22, 48.32561117432184, 46.04709243774414, 51.8484251691426, 45.8774131497029, 212968277),
    @(23, 38.91063368940384, 41.01258087158203, 43.73854191057829, 37.31775382246912, 212968277),
    @(24, 40.22891917780825, 37.60909271240234, 42.26656409933605, 36.20353196964999, 212968277),
    @(25, 42.00091843861417, 40.31481170654297, 44.68182839988685, 39.55606239071358, 212968277),
    @(26, 53.06029015741621, 56.32278823852539, 65.65028874834994, 49.44002639349635, 212968277),
    @(27, 66.39796070487871, 79.78755187988281, 83.46582898942557, 66.00263323541735, 212968277),
    @(28, 76.2818797263953, 76.10871124267578, 78.55042385646477, 68.28137555018625, 212968277),
    @(29, 73.03608249391461, 77.74977874755859, 78.55285109884113, 68.26128305397191, 212968277),
    @(30, 99.3259585367214, 94.15886688232422, 103.4033927751003, 80.74902563820369, 212968277),
    @(31, 80.11333617339662, 72.63231658935547, 81.27312055279506, 69.89664930620528, 212968277),
    @(32, 63.48782397067654, 71.50762176513672, 74.90507892900976, 60.22447166260122, 212968277),
    @(33, 48.44201907237223, 44.90630340576172, 52.44795822668351, 44.83396318817324, 212968277),
    @(34, 48.4469940503439, 62.14491271972656, 65.09720983424245, 46.70495464167461, 212968277),
    @(35, 61.2032349567249, 54.49869918823242, 62.89096024374754, 50.37211672987199, 212968277),
    @(36, 58.12674683097779, 59.53303909301758, 60.98621062490095, 53.73911468730277, 212968277),
    @(38, 80.93645903584445, 81.87291717529297, 87.95031915185052, 75.86240454509661, 212968277),
    @(39, 89.73660419246195, 82.76190948486328, 90.88299272372689, 79.29382234036787, 212968277),
    @(40, 99.88533207244586, 99.08042907714844, 110.1357003202436, 95.65717567631934, 212968277),
    @(41, 106.678686533637, 97.98975372314452, 112.5852153056132, 96.56437814469898, 212968277),
    @(42, 85.7743320136673, 94.77370452880859, 110.0677205656978, 84.37771011617207, 212968277),
    @(43, 83.8887094637239, 90.2646484375, 90.74061608380831, 62.65871786546948, 212968277),
    @(44, 143.1886625997294, 156.5166473388672, 156.955267231608, 137.8654452705798, 212968277)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("D$rowNum").Value = $r[1]
    $ws.Range("E$rowNum").Value = $r[2]
    $ws.Range("F$rowNum").Value = $r[3]
    $ws.Range("G$rowNum").Value = $r[4]
    $ws.Range("H$rowNum").Value = $r[5]
    $ws.Range("I$rowNum").Value = "STX"
}
